$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Copy the date-formatted cell above onto the two new date cells first so
# they pick up the existing date number-format style (reused, no new xf).
$ws.Range("B32").Copy($ws.Range("B33"))
$ws.Range("B32").Copy($ws.Range("B34"))

# Row 33: Federico Speroni, 2017-05-08 (serial 42863), 3 hours
$ws.Range("A33").Value = "Federico Speroni"
$ws.Range("B33").Value = 42863
$ws.Range("C33").Value = 3
$ws.Range("D33").Value = "Sprint 1 - Diagrama de Clases"
$ws.Range("E33").Value = "Idea general de Diagrama de clases, creación de Diagrama de clases para Sprint1"

# Row 34: Federico Speroni, 2017-05-09 (serial 42864), 6 hours
$ws.Range("A34").Value = "Federico Speroni"
$ws.Range("B34").Value = 42864
$ws.Range("C34").Value = 6
$ws.Range("D34").Value = "Sprint 1 - Back-end"
$ws.Range("E34").Value = "Armado de las capas y estructura del back-end, investigación de servicios web, realización de la interfas del servicio."

# Update selection / view to match the final state
$ws.Range("E34").Select()
